$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -4
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = -2
